$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 152.8
$ws.Range("I9").Value = 134.36363
$ws.Range("K9").Value = 134.36363
$ws.Range("M9").Value = 34.63637
$ws.Range("H19").Value = 1056.1724
$ws.Range("I19").Value = 525.9231
$ws.Range("J19").Value = 1487
$ws.Range("K19").Value = 525.9231
$ws.Range("L19").Value = 1487
$ws.Range("M19").Value = -350.9231
$ws.Range("N19").Value = -1837
$ws.Range("H32").Value = 5071.304
$ws.Range("I32").Value = 6332.3335
$ws.Range("J32").Value = 4882.15
$ws.Range("K32").Value = 6332.3335
$ws.Range("L32").Value = 4882.15
$ws.Range("M32").Value = -6006.3335
$ws.Range("N32").Value = -5534.15
$ws.Range("H43").Value = 1047.4445
$ws.Range("I43").Value = 1040.7142
$ws.Range("K43").Value = 1040.7142
$ws.Range("M43").Value = -971.7141999999999
$ws.Range("H55").Value = 180.14285
$ws.Range("I55").Value = 151.8
$ws.Range("J55").Value = 195.88889
$ws.Range("K55").Value = 151.8
$ws.Range("L55").Value = 195.88889
$ws.Range("M55").Value = 62.19999999999999
$ws.Range("N55").Value = -623.8888899999999
$ws.Range("H62").Value = 3307.923
$ws.Range("I62").Value = 3307.923
$ws.Range("K62").Value = 3307.923
$ws.Range("M62").Value = -2683.923
$ws.Range("H65").Value = 3307.923
$ws.Range("I65").Value = 3307.923
$ws.Range("K65").Value = 16539.615
$ws.Range("M65").Value = -13419.615
$ws.Range("H80").Value = 882.2857
$ws.Range("I80").Value = 476.85715
$ws.Range("J80").Value = 1287.7142
$ws.Range("K80").Value = 1430.57145
$ws.Range("L80").Value = 3863.1426
$ws.Range("M80").Value = -432.5714499999999
$ws.Range("N80").Value = -5859.142599999999
$ws.Range("H83").Value = 882.2857
$ws.Range("I83").Value = 476.85715
$ws.Range("J83").Value = 1287.7142
$ws.Range("K83").Value = 4291.71435
$ws.Range("L83").Value = 11589.4278
$ws.Range("M83").Value = 700.2856499999998
$ws.Range("N83").Value = -21573.4278
$ws.Range("H132").Value = 37040092
$ws.Range("I132").Value = 40003060
$ws.Range("J132").Value = 2998.5
$ws.Range("K132").Value = 120009180
$ws.Range("L132").Value = 8995.5
$ws.Range("M132").Value = -120006650
$ws.Range("N132").Value = -14055.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 933.4167
$ws.Range("J4").Value = 1002
$ws.Range("L4").Value = 1002
$ws.Range("N4").Value = -1234
$ws.Range("H32").Value = 5374.717
$ws.Range("I32").Value = 3110.7805
$ws.Range("K32").Value = 3110.7805
$ws.Range("M32").Value = -2823.7805
$ws.Range("H61").Value = 2982.5
$ws.Range("I61").Value = 2482.5
$ws.Range("K61").Value = 2482.5
$ws.Range("M61").Value = -2270.5
$ws.Range("H122").Value = 445412.97
$ws.Range("I122").Value = 1868.5883
$ws.Range("J122").Value = 1605452.1
$ws.Range("K122").Value = 5605.7649
$ws.Range("L122").Value = 4816356.300000001
$ws.Range("M122").Value = -3155.7649
$ws.Range("N122").Value = -4821256.300000001
$ws.Range("H132").Value = 2643.1667
$ws.Range("I132").Value = 2162.55
$ws.Range("K132").Value = 6487.650000000001
$ws.Range("M132").Value = -3957.650000000001
$ws.Range("H136").Value = 2982.5
$ws.Range("I136").Value = 2482.5
$ws.Range("K136").Value = 7447.5
$ws.Range("M136").Value = -4897.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 32419.8
$ws.Range("J23").Value = 32419.8
$ws.Range("L23").Value = 32419.8
$ws.Range("N23").Value = -32985.8
$ws.Range("H25").Value = 1750
$ws.Range("I25").Value = 1750
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 1750
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -1515
$ws.Range("N25").ClearContents()
$ws.Range("H29").Value = 155857.14
$ws.Range("I29").Value = 148500
$ws.Range("K29").Value = 148500
$ws.Range("M29").Value = -148211
$ws.Range("H107").Value = 4469625
$ws.Range("J107").Value = 7839
$ws.Range("L107").Value = 7839
$ws.Range("N107").Value = -11679
$ws.Range("H109").Value = 59994
$ws.Range("J109").Value = 59994
$ws.Range("L109").Value = 59994
$ws.Range("N109").Value = -62768
$ws.Range("H134").Value = 3947.818
$ws.Range("I134").Value = 1886.7916
$ws.Range("J134").Value = 9443.888999999999
$ws.Range("K134").Value = 5660.3748
$ws.Range("L134").Value = 28331.667
$ws.Range("M134").Value = -3125.3748
$ws.Range("N134").Value = -33401.667

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 3768.077
$ws.Range("I5").Value = 1984.4286
$ws.Range("J5").Value = 5849
$ws.Range("K5").Value = 1984.4286
$ws.Range("L5").Value = 5849
$ws.Range("M5").Value = -1872.4286
$ws.Range("N5").Value = -6073
$ws.Range("H14").Value = 660
$ws.Range("I14").Value = 660
$ws.Range("K14").Value = 660
$ws.Range("M14").Value = -490
$ws.Range("I15").Value = 1000
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 1000
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -830
$ws.Range("N15").ClearContents()
$ws.Range("H25").Value = 2000
$ws.Range("J25").Value = 2000
$ws.Range("L25").Value = 2000
$ws.Range("N25").Value = -2348
$ws.Range("H31").Value = 15956.539
$ws.Range("I31").Value = 2169.6667
$ws.Range("J31").Value = 19346.754
$ws.Range("K31").Value = 2169.6667
$ws.Range("L31").Value = 19346.754
$ws.Range("M31").Value = -1874.6667
$ws.Range("N31").Value = -19936.754
$ws.Range("H34").Value = 15956.539
$ws.Range("I34").Value = 2169.6667
$ws.Range("J34").Value = 19346.754
$ws.Range("K34").Value = 2169.6667
$ws.Range("L34").Value = 19346.754
$ws.Range("M34").Value = -1967.6667
$ws.Range("N34").Value = -19750.754
$ws.Range("H105").Value = 1969.75
$ws.Range("I105").Value = 1793
$ws.Range("K105").Value = 1793
$ws.Range("M105").Value = -46
$ws.Range("H134").Value = 2634.0386
$ws.Range("I134").Value = 1826.4445
$ws.Range("J134").Value = 4451.125
$ws.Range("K134").Value = 5479.333500000001
$ws.Range("L134").Value = 13353.375
$ws.Range("M134").Value = -2944.333500000001
$ws.Range("N134").Value = -18423.375

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 74290.336
$ws.Range("J12").Value = 292
$ws.Range("L12").Value = 876
$ws.Range("N12").Value = -1222
$ws.Range("H86").Value = 418.6
$ws.Range("J86").Value = 423.5
$ws.Range("L86").Value = 1270.5
$ws.Range("N86").Value = -3642.5
$ws.Range("H89").Value = 418.6
$ws.Range("J89").Value = 423.5
$ws.Range("L89").Value = 3811.5
$ws.Range("N89").Value = -15667.5
$ws.Range("H105").Value = 5000
$ws.Range("J105").Value = 5000
$ws.Range("L105").Value = 15000
$ws.Range("N105").Value = -20242
$ws.Range("H122").Value = 839.7143
$ws.Range("I122").Value = 555.1818
$ws.Range("J122").Value = 1883
$ws.Range("K122").Value = 4996.6362
$ws.Range("L122").Value = 16947
$ws.Range("M122").Value = -2546.6362
$ws.Range("N122").Value = -21847

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 6756.8
$ws.Range("I2").Value = 116.454544
$ws.Range("K2").Value = 116.454544
$ws.Range("M2").Value = -3.454543999999999
$ws.Range("H62").Value = 1234567
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 1234567
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H70").Value = 13341000
$ws.Range("I70").Value = 20008300
$ws.Range("J70").Value = 6399.6
$ws.Range("K70").Value = 20008300
$ws.Range("L70").Value = 6399.6
$ws.Range("M70").Value = -20008030
$ws.Range("N70").Value = -6939.6
$ws.Range("H73").Value = 13341000
$ws.Range("I73").Value = 20008300
$ws.Range("J73").Value = 6399.6
$ws.Range("K73").Value = 20008300
$ws.Range("L73").Value = 6399.6
$ws.Range("M73").Value = -20007364
$ws.Range("N73").Value = -8271.6
$ws.Range("H80").Value = 1438115.1
$ws.Range("I80").Value = 2712294
$ws.Range("K80").Value = 2712294
$ws.Range("M80").Value = -2711296
$ws.Range("H83").Value = 1438115.1
$ws.Range("I83").Value = 2712294
$ws.Range("K83").Value = 13561470
$ws.Range("M83").Value = -13556478

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1605
$ws.Range("I9").Value = 760
$ws.Range("K9").Value = 760
$ws.Range("M9").Value = -536
$ws.Range("H30").Value = 142859860
$ws.Range("I30").Value = 142859860
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 142859860
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -142859752
$ws.Range("N30").ClearContents()
$ws.Range("H82").Value = 4275395
$ws.Range("J82").Value = 1310.1666
$ws.Range("L82").Value = 1310.1666
$ws.Range("N82").Value = -2032.1666
$ws.Range("H85").Value = 4275395
$ws.Range("J85").Value = 1310.1666
$ws.Range("L85").Value = 1310.1666
$ws.Range("N85").Value = -3806.1666
$ws.Range("H136").Value = 81492.46000000001
$ws.Range("I136").Value = 129050.5
$ws.Range("K136").Value = 387151.5
$ws.Range("M136").Value = -384601.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 91001820
$ws.Range("I132").Value = 100013120
$ws.Range("K132").Value = 300039360
$ws.Range("M132").Value = -300036830
$ws.Range("H136").Value = 3835.1765
$ws.Range("I136").Value = 3215.2307
$ws.Range("K136").Value = 9645.6921
$ws.Range("M136").Value = -7095.6921
